$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the new rows with the "Good" cell style (same as existing data rows)
$ws.Range("A44:AO46").Style = "Good"

# Populate new rows 44-46 with the additional effect-size records
    # Row 44
    $ws.Range("A44").Value = "Good"
    $ws.Range("B44").Value = "TEMP638"
    $ws.Range("C44").Value = 582
    $ws.Range("D44").Value = 10
    $ws.Range("E44").Value = "20B"
    $ws.Range("F44").Value = "Gharibi, M. R.et al."
    $ws.Range("G44").Value = 2016
    $ws.Range("H44").Value = "Life Cycle of the Fairy Shrimp, Phallocryptus spinosa Milne Edwards, 1840 (Crustacea: Anostraca) at Different Temperatures"
    $ws.Range("I44").Value = "NA"
    $ws.Range("J44").Value = "Phallocryptus spinosa"
    $ws.Range("K44").Value = "Thamnocephalidae"
    $ws.Range("L44").Value = "Crustacea"
    $ws.Range("M44").Value = "Arthropoda"
    $ws.Range("N44").Value = "Both"
    $ws.Range("O44").Value = "Internal"
    $ws.Range("P44").Value = "No"
    $ws.Range("Q44").Value = "Lab"
    $ws.Range("R44").Value = "Wild-caught"
    $ws.Range("S44").Value = "Experimental"
    $ws.Range("T44").Value = "Aquatic"
    $ws.Range("U44").Value = "Iran"
    $ws.Range("V44").Value = "Asia"
    $ws.Range("W44").Value = "NA"
    $ws.Range("X44").Value = "More than 5 days"
    $ws.Range("Y44").Value = "Mix"
    $ws.Range("Z44").Value = "No"
    $ws.Range("AA44").Value = "Constant"
    $ws.Range("AB44").Value = "No"
    $ws.Range("AC44").Value = "TEMP638_1"
    $ws.Range("AD44").Value = "TEMP638_A"
    $ws.Range("AE44").Value = "Reproduction"
    $ws.Range("AF44").Value = "number of clutch per female"
    $ws.Range("AG44").Value = "Female"
    $ws.Range("AH44").Value = "Female"
    $ws.Range("AI44").Value = 25
    $ws.Range("AJ44").Value = -2.2975181029999998
    $ws.Range("AK44").Value = 0.20747795999999999
    $ws.Range("AL44").Value = -10
    $ws.Range("AM44").Value = 25
    $ws.Range("AN44").Value = 15
    $ws.Range("AO44").Value = "Cool"

    # Row 45
    $ws.Range("A45").Value = "Good"
    $ws.Range("B45").Value = "TEMP638"
    $ws.Range("C45").Value = 582
    $ws.Range("D45").Value = 10
    $ws.Range("E45").Value = "20B"
    $ws.Range("F45").Value = "Gharibi, M. R.et al."
    $ws.Range("G45").Value = 2016
    $ws.Range("H45").Value = "Life Cycle of the Fairy Shrimp, Phallocryptus spinosa Milne Edwards, 1840 (Crustacea: Anostraca) at Different Temperatures"
    $ws.Range("I45").Value = "NA"
    $ws.Range("J45").Value = "Phallocryptus spinosa"
    $ws.Range("K45").Value = "Thamnocephalidae"
    $ws.Range("L45").Value = "Crustacea"
    $ws.Range("M45").Value = "Arthropoda"
    $ws.Range("N45").Value = "Both"
    $ws.Range("O45").Value = "Internal"
    $ws.Range("P45").Value = "No"
    $ws.Range("Q45").Value = "Lab"
    $ws.Range("R45").Value = "Wild-caught"
    $ws.Range("S45").Value = "Experimental"
    $ws.Range("T45").Value = "Aquatic"
    $ws.Range("U45").Value = "Iran"
    $ws.Range("V45").Value = "Asia"
    $ws.Range("W45").Value = "NA"
    $ws.Range("X45").Value = "More than 5 days"
    $ws.Range("Y45").Value = "Mix"
    $ws.Range("Z45").Value = "No"
    $ws.Range("AA45").Value = "Constant"
    $ws.Range("AB45").Value = "No"
    $ws.Range("AC45").Value = "TEMP638_1"
    $ws.Range("AD45").Value = "TEMP638_A"
    $ws.Range("AE45").Value = "Reproduction"
    $ws.Range("AF45").Value = "number of clutch per female"
    $ws.Range("AG45").Value = "Female"
    $ws.Range("AH45").Value = "Female"
    $ws.Range("AI45").Value = 25
    $ws.Range("AJ45").Value = 0.17083727600000001
    $ws.Range("AK45").Value = 0.125456022
    $ws.Range("AL45").Value = -5
    $ws.Range("AM45").Value = 25
    $ws.Range("AN45").Value = 20
    $ws.Range("AO45").Value = "Cool"

    # Row 46
    $ws.Range("A46").Value = "Good"
    $ws.Range("B46").Value = "TEMP638"
    $ws.Range("C46").Value = 582
    $ws.Range("D46").Value = 10
    $ws.Range("E46").Value = "20B"
    $ws.Range("F46").Value = "Gharibi, M. R.et al."
    $ws.Range("G46").Value = 2016
    $ws.Range("H46").Value = "Life Cycle of the Fairy Shrimp, Phallocryptus spinosa Milne Edwards, 1840 (Crustacea: Anostraca) at Different Temperatures"
    $ws.Range("I46").Value = "NA"
    $ws.Range("J46").Value = "Phallocryptus spinosa"
    $ws.Range("K46").Value = "Thamnocephalidae"
    $ws.Range("L46").Value = "Crustacea"
    $ws.Range("M46").Value = "Arthropoda"
    $ws.Range("N46").Value = "Both"
    $ws.Range("O46").Value = "Internal"
    $ws.Range("P46").Value = "No"
    $ws.Range("Q46").Value = "Lab"
    $ws.Range("R46").Value = "Wild-caught"
    $ws.Range("S46").Value = "Experimental"
    $ws.Range("T46").Value = "Aquatic"
    $ws.Range("U46").Value = "Iran"
    $ws.Range("V46").Value = "Asia"
    $ws.Range("W46").Value = "NA"
    $ws.Range("X46").Value = "More than 5 days"
    $ws.Range("Y46").Value = "Mix"
    $ws.Range("Z46").Value = "No"
    $ws.Range("AA46").Value = "Constant"
    $ws.Range("AB46").Value = "No"
    $ws.Range("AC46").Value = "TEMP638_1"
    $ws.Range("AD46").Value = "TEMP638_A"
    $ws.Range("AE46").Value = "Reproduction"
    $ws.Range("AF46").Value = "number of clutch per female"
    $ws.Range("AG46").Value = "Female"
    $ws.Range("AH46").Value = "Female"
    $ws.Range("AI46").Value = 25
    $ws.Range("AJ46").Value = 0
    $ws.Range("AK46").Value = 0
    $ws.Range("AL46").Value = 0
    $ws.Range("AM46").Value = 25
    $ws.Range("AN46").Value = 25
    $ws.Range("AO46").Value = "Reference"

# Center-align the "es" column values for the new rows (new cell style with horizontal center alignment)
$ws.Range("AI44:AI46").HorizontalAlignment = -4108

# Update selection to match the edited location
$ws.Range("A45").Select()
